$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 686.8946999999999
$ws.Range("J19").Value = 729.5333000000001
$ws.Range("L19").Value = 729.5333000000001
$ws.Range("N19").Value = -1079.5333

$ws.Range("H41").Value = 560.5454999999999
$ws.Range("I41").Value = 135
$ws.Range("J41").Value = 915.1667
$ws.Range("K41").Value = 135
$ws.Range("L41").Value = 915.1667
$ws.Range("M41").Value = 305
$ws.Range("N41").Value = -1795.1667

$ws.Range("H74").Value = 3866
$ws.Range("I74").Value = 3500
$ws.Range("J74").Value = 4075.1428
$ws.Range("K74").Value = 3500
$ws.Range("L74").Value = 4075.1428
$ws.Range("M74").Value = -2564
$ws.Range("N74").Value = -5947.1428

$ws.Range("H77").Value = 3866
$ws.Range("I77").Value = 3500
$ws.Range("J77").Value = 4075.1428
$ws.Range("K77").Value = 17500
$ws.Range("L77").Value = 20375.714
$ws.Range("M77").Value = -12820
$ws.Range("N77").Value = -29735.714

$ws.Range("H86").Value = 3200.9092
$ws.Range("I86").Value = 1241.2
$ws.Range("J86").Value = 4834
$ws.Range("K86").Value = 1241.2
$ws.Range("L86").Value = 4834
$ws.Range("M86").Value = -118.2
$ws.Range("N86").Value = -7080

$ws.Range("H89").Value = 3200.9092
$ws.Range("I89").Value = 1241.2
$ws.Range("J89").Value = 4834
$ws.Range("K89").Value = 6206
$ws.Range("L89").Value = 24170
$ws.Range("M89").Value = -590
$ws.Range("N89").Value = -35402

$ws.Range("H125").Value = 1179.1666
$ws.Range("I125").Value = 979.6667
$ws.Range("J125").Value = 1378.6666
$ws.Range("K125").Value = 8817.0003
$ws.Range("L125").Value = 12407.9994
$ws.Range("M125").Value = -6357.0003
$ws.Range("N125").Value = -17327.9994

$ws.Range("H137").Value = 1030.2051
$ws.Range("I137").Value = 1030.7273
$ws.Range("J137").Value = 1029.5294
$ws.Range("K137").Value = 3092.1819
$ws.Range("L137").Value = 3088.5882
$ws.Range("M137").Value = -542.1819
$ws.Range("N137").Value = -8188.5882

$ws.Range("H138").Value = 3101.239
$ws.Range("I138").Value = 1118.3914
$ws.Range("J138").Value = 5084.087
$ws.Range("K138").Value = 3355.1742
$ws.Range("L138").Value = 15252.261
$ws.Range("M138").Value = 1784.8258
$ws.Range("N138").Value = -25532.261

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 446921.7
$ws.Range("I32").Value = 4270.5474
$ws.Range("J32").Value = 2401964.2
$ws.Range("K32").Value = 4270.5474
$ws.Range("L32").Value = 2401964.2
$ws.Range("M32").Value = -3983.5474
$ws.Range("N32").Value = -2402538.2

$ws.Range("H61").Value = 2320.7585
$ws.Range("I61").Value = 2271.5454
$ws.Range("J61").Value = 2475.4285
$ws.Range("K61").Value = 2271.5454
$ws.Range("L61").Value = 2475.4285
$ws.Range("M61").Value = -2059.5454
$ws.Range("N61").Value = -2899.4285

$ws.Range("H110").Value = 728.8461
$ws.Range("I110").Value = 731.2
$ws.Range("J110").Value = 721
$ws.Range("K110").Value = 731.2
$ws.Range("L110").Value = 721
$ws.Range("M110").Value = 1313.8
$ws.Range("N110").Value = -4811

$ws.Range("H136").Value = 2320.7585
$ws.Range("I136").Value = 2271.5454
$ws.Range("J136").Value = 2475.4285
$ws.Range("K136").Value = 6814.6362
$ws.Range("L136").Value = 7426.2855
$ws.Range("M136").Value = -4264.6362
$ws.Range("N136").Value = -12526.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1672.1666
$ws.Range("I134").Value = 1569.4
$ws.Range("J134").Value = 2802.6
$ws.Range("K134").Value = 4708.200000000001
$ws.Range("L134").Value = 8407.799999999999
$ws.Range("M134").Value = -2173.200000000001
$ws.Range("N134").Value = -13477.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 69.416664
$ws.Range("I7").Value = 47.166668
$ws.Range("K7").Value = 47.166668
$ws.Range("M7").Value = 65.833332

$ws.Range("H16").Value = 1789.409
$ws.Range("I16").Value = 1409.2778
$ws.Range("J16").Value = 3500
$ws.Range("K16").Value = 1409.2778
$ws.Range("L16").Value = 3500
$ws.Range("M16").Value = -1122.2778
$ws.Range("N16").Value = -4074

$ws.Range("H99").Value = 396968.2
$ws.Range("I99").Value = 448978.9
$ws.Range("J99").Value = 1686.8
$ws.Range("K99").Value = 448978.9
$ws.Range("L99").Value = 1686.8
$ws.Range("M99").Value = -447480.9
$ws.Range("N99").Value = -4682.8

$ws.Range("H105").Value = 1567.6154
$ws.Range("I105").Value = 1015.6
$ws.Range("J105").Value = 1912.625
$ws.Range("K105").Value = 1015.6
$ws.Range("L105").Value = 1912.625
$ws.Range("M105").Value = 731.4
$ws.Range("N105").Value = -5406.625

$ws.Range("H113").Value = 1789.409
$ws.Range("I113").Value = 1409.2778
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 1409.2778
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = 760.7221999999999
$ws.Range("N113").Value = -7840

$ws.Range("H126").Value = 396968.2
$ws.Range("I126").Value = 448978.9
$ws.Range("J126").Value = 1686.8
$ws.Range("K126").Value = 1346936.7
$ws.Range("L126").Value = 5060.4
$ws.Range("M126").Value = -1344466.7
$ws.Range("N126").Value = -10000.4

$ws.Range("H134").Value = 150019000
$ws.Range("I134").Value = 240001200
$ws.Range("J134").Value = 48671.332
$ws.Range("K134").Value = 720003600
$ws.Range("L134").Value = 146013.996
$ws.Range("M134").Value = -720001065
$ws.Range("N134").Value = -151083.996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 28145
$ws.Range("I46").Value = 1083.3334
$ws.Range("K46").Value = 3250.0002
$ws.Range("M46").Value = -3159.0002

$ws.Range("H68").Value = 954.8214
$ws.Range("I68").Value = 775
$ws.Range("J68").Value = 984.7917
$ws.Range("K68").Value = 2325
$ws.Range("L68").Value = 2954.3751
$ws.Range("M68").Value = -1514
$ws.Range("N68").Value = -4576.3751

$ws.Range("H70").Value = 2000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H71").Value = 954.8214
$ws.Range("I71").Value = 775
$ws.Range("J71").Value = 984.7917
$ws.Range("K71").Value = 6975
$ws.Range("L71").Value = 8863.1253
$ws.Range("M71").Value = -2919
$ws.Range("N71").Value = -16975.1253

$ws.Range("H73").Value = 2000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H122").Value = 420.4737
$ws.Range("J122").Value = 563.8
$ws.Range("L122").Value = 5074.2
$ws.Range("N122").Value = -9974.200000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 3000
$ws.Range("J21").Value = 3000
$ws.Range("L21").Value = 3000
$ws.Range("N21").Value = -3348

$ws.Range("H122").Value = 100005980
$ws.Range("I122").Value = 7475.25
$ws.Range("J122").Value = 500000000
$ws.Range("K122").Value = 22425.75
$ws.Range("L122").Value = 1500000000
$ws.Range("M122").Value = -19975.75
$ws.Range("N122").Value = -1500004900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7237
$ws.Range("I81").Value = 2412.8572
$ws.Range("J81").Value = 12865.167
$ws.Range("K81").Value = 4825.7144
$ws.Range("L81").Value = 25730.334
$ws.Range("M81").Value = -3764.7144
$ws.Range("N81").Value = -27852.334

$ws.Range("H84").Value = 7237
$ws.Range("I84").Value = 2412.8572
$ws.Range("J84").Value = 12865.167
$ws.Range("K84").Value = 24128.572
$ws.Range("L84").Value = 128651.67
$ws.Range("M84").Value = -18824.572
$ws.Range("N84").Value = -139259.67

$ws.Range("H132").Value = 31747720
$ws.Range("I132").Value = 44446040
$ws.Range("J132").Value = 1917.1111
$ws.Range("K132").Value = 133338120
$ws.Range("L132").Value = 5751.3333
$ws.Range("M132").Value = -133335590
$ws.Range("N132").Value = -10811.3333
